# Fix client payment report: insert a new "CLIENT" column before the
# existing "LINE OF BUSINESS" column (i.e. before column B), shifting
# all the other header columns one place to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B; this shifts B:G to C:H, preserving their
# widths, values and styles.
$ws.Columns.Item(2).Insert()

# Give the new column B the same header formatting (style/border/fill)
# as the rest of the header row by copying column A's header cell,
# then overwrite the text with "CLIENT".
$ws.Cells.Item(1, 1).Copy()
$ws.Cells.Item(1, 2).PasteSpecial()
$ws.Cells.Item(1, 2).Value = "CLIENT"

# Set the new column's width as specified by the report layout
# (27.5 characters once Excel's internal width conversion is applied).
$ws.Columns.Item(2).ColumnWidth = 26.67
